# Updates cryptos list values (Price / Volume(1h) columns) and fixes the
# order of three coin rows (45-47), matching the "Thu Dec 14 04:39:07 UTC
# 2023" GitHub Actions refresh.
#
# Note: several "Price" values look like plain numbers (e.g. "248.69").
# The original cells store them as literal text, so a bare
# `Range.Value = "248.69"` would be auto-coerced to a numeric cell by
# Excel. To keep them as text we assign with a leading apostrophe (the
# same trick used when typing numbers-as-text directly into Excel) and
# then reset the cell style to "Normal" so no stray number-format/
# quote-prefix styling is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.725.74'
$ws.Range('E2').Value = '  +4.65%  '
$ws.Range('D3').Value = '2.250.75'
$ws.Range('E3').Value = '  +3.95%  '
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').Value = "'248.69"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.81%  '
$ws.Range('D6').Value = "'0.634"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +4.52%  '
$ws.Range('D7').Value = "'70.33"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +6.48%  '
$ws.Range('E8').Value = '  +0.08%  '
$ws.Range('D9').Value = "'0.651"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +15.22%  '
$ws.Range('E10').Value = '  +10.42%  '
$ws.Range('D11').Value = "'59.47"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +2.28%  '
$ws.Range('E12').Value = '  +4.69%  '
$ws.Range('D13').Value = "'7.47"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +8.67%  '
$ws.Range('D14').Value = "'0.105"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.51%  '
$ws.Range('D15').Value = '2.580.19'
$ws.Range('E15').Value = '  +3.98%  '
$ws.Range('D16').Value = "'14.82"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.14%  '
$ws.Range('D17').Value = "'0.877"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('D18').Value = '2.268.09'
$ws.Range('E18').Value = '  +4.70%  '
$ws.Range('D19').Value = '42.682.88'
$ws.Range('E19').Value = '  +4.64%  '
$ws.Range('D20').Value = '0.0₃0989'
$ws.Range('E20').Value = '  +5.80%  '
$ws.Range('E21').Value = '  +3.54%  '
$ws.Range('D22').Value = "'72.95"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.58%  '
$ws.Range('D23').Value = "'235.28"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.12%  '
$ws.Range('D24').Value = "'2.04"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -2.88%  '
$ws.Range('D25').Value = "'3.93"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +6.70%  '
$ws.Range('D26').Value = "'11.48"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.30%  '
$ws.Range('E27').Value = '  -0.04%  '
$ws.Range('D28').Value = "'2.42"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.05%  '
$ws.Range('E29').Value = '  -1.46%  '
$ws.Range('D30').Value = "'2.20"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +10.79%  '
$ws.Range('D31').Value = "'167.65"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.30%  '
$ws.Range('D32').Value = "'20.83"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.72%  '
$ws.Range('D33').Value = "'6.43"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +15.76%  '
$ws.Range('D34').Value = "'0.126"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +5.56%  '
$ws.Range('D35').Value = "'0.0804"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +9.48%  '
$ws.Range('D36').Value = "'31.34"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +26.13%  '
$ws.Range('E37').Value = '  +3.75%  '
$ws.Range('D38').Value = "'4.43"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +12.20%  '
$ws.Range('D39').Value = "'4.69"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.73%  '
$ws.Range('E40').Value = '  +8.42%  '
$ws.Range('D41').Value = "'2.29"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.90%  '
$ws.Range('D42').Value = "'12.44"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +7.48%  '
$ws.Range('D43').Value = "'5.80"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +6.70%  '
$ws.Range('D44').Value = "'61.96"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.10%  '
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = "'0.202"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.95%  '
$ws.Range('B46').Value = 'FraxShare'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D46').Value = "'8.99"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +6.45%  '
$ws.Range('B47').Value = 'FTXToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range('D47').Value = "'4.89"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('E48').Value = '  +2.88%  '
$ws.Range('E49').Value = '  -0.13%  '
$ws.Range('D50').Value = "'1.16"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.26%  '
$ws.Range('E51').Value = '  +4.20%  '
